$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "40+25=65"
$t.Cell(1,2).Range.Text = "68-37=31"
$t.Cell(1,3).Range.Text = "31+53=84"
$t.Cell(1,4).Range.Text = "49+20=69"
$t.Cell(1,5).Range.Text = "88+4=92"
$t.Cell(2,1).Range.Text = "99-14=85"
$t.Cell(2,2).Range.Text = "33+5=38"
$t.Cell(2,3).Range.Text = "90-30=60"
$t.Cell(2,4).Range.Text = "9+89=98"
$t.Cell(2,5).Range.Text = "94-41=53"
$t.Cell(3,1).Range.Text = "13+0=13"
$t.Cell(3,2).Range.Text = "21+46=67"
$t.Cell(3,3).Range.Text = "56-37=19"
$t.Cell(3,4).Range.Text = "85-71=14"
$t.Cell(3,5).Range.Text = "72+6=78"
$t.Cell(4,1).Range.Text = "23+30=53"
$t.Cell(4,2).Range.Text = "2+49=51"
$t.Cell(4,3).Range.Text = "78-0=78"
$t.Cell(4,4).Range.Text = "28+12=40"
$t.Cell(4,5).Range.Text = "64+30=94"
$t.Cell(5,1).Range.Text = "24-2=22"
$t.Cell(5,2).Range.Text = "36+58=94"
$t.Cell(5,3).Range.Text = "0+60=60"
$t.Cell(5,4).Range.Text = "87+3=90"
$t.Cell(5,5).Range.Text = "29+26=55"
$t.Cell(6,1).Range.Text = "87+11=98"
$t.Cell(6,2).Range.Text = "97-29=68"
$t.Cell(6,3).Range.Text = "91-23=68"
$t.Cell(6,4).Range.Text = "54-47=7"
$t.Cell(6,5).Range.Text = "66-38=28"
$t.Cell(7,1).Range.Text = "92-23=69"
$t.Cell(7,2).Range.Text = "24+18=42"
$t.Cell(7,3).Range.Text = "13+9=22"
$t.Cell(7,4).Range.Text = "4+50=54"
$t.Cell(7,5).Range.Text = "49-11=38"
$t.Cell(8,1).Range.Text = "65-34=31"
$t.Cell(8,2).Range.Text = "15+36=51"
$t.Cell(8,3).Range.Text = "1+58=59"
$t.Cell(8,4).Range.Text = "17+7=24"
$t.Cell(8,5).Range.Text = "9-0=9"
$t.Cell(9,1).Range.Text = "47-39=8"
$t.Cell(9,2).Range.Text = "0+75=75"
$t.Cell(9,3).Range.Text = "97-9=88"
$t.Cell(9,4).Range.Text = "34+56=90"
$t.Cell(9,5).Range.Text = "20+21=41"
$t.Cell(10,1).Range.Text = "97-20=77"
$t.Cell(10,2).Range.Text = "1+88=89"
$t.Cell(10,3).Range.Text = "32-18=14"
$t.Cell(10,4).Range.Text = "19+71=90"
$t.Cell(10,5).Range.Text = "98-82=16"
$t.Cell(11,1).Range.Text = "56+5=61"
$t.Cell(11,2).Range.Text = "43-36=7"
$t.Cell(11,3).Range.Text = "41+30=71"
$t.Cell(11,4).Range.Text = "37+6=43"
$t.Cell(11,5).Range.Text = "28-27=1"
$t.Cell(12,1).Range.Text = "50+37=87"
$t.Cell(12,2).Range.Text = "73-46=27"
$t.Cell(12,3).Range.Text = "20+15=35"
$t.Cell(12,4).Range.Text = "28+16=44"
$t.Cell(12,5).Range.Text = "54+42=96"
$t.Cell(13,1).Range.Text = "23+52=75"
$t.Cell(13,2).Range.Text = "65+24=89"
$t.Cell(13,3).Range.Text = "13+44=57"
$t.Cell(13,4).Range.Text = "51+9=60"
$t.Cell(13,5).Range.Text = "84-56=28"
$t.Cell(14,1).Range.Text = "3+44=47"
$t.Cell(14,2).Range.Text = "30+27=57"
$t.Cell(14,3).Range.Text = "28+70=98"
$t.Cell(14,4).Range.Text = "59-18=41"
$t.Cell(14,5).Range.Text = "28+47=75"
$t.Cell(15,1).Range.Text = "77-53=24"
$t.Cell(15,2).Range.Text = "13+24=37"
$t.Cell(15,3).Range.Text = "37+15=52"
$t.Cell(15,4).Range.Text = "59-1=58"
$t.Cell(15,5).Range.Text = "84-21=63"
$t.Cell(16,1).Range.Text = "36+61=97"
$t.Cell(16,2).Range.Text = "72+0=72"
$t.Cell(16,3).Range.Text = "86-25=61"
$t.Cell(16,4).Range.Text = "67-25=42"
$t.Cell(16,5).Range.Text = "35-27=8"
$t.Cell(17,1).Range.Text = "78-51=27"
$t.Cell(17,2).Range.Text = "2+36=38"
$t.Cell(17,3).Range.Text = "96-26=70"
$t.Cell(17,4).Range.Text = "77-29=48"
$t.Cell(17,5).Range.Text = "36-19=17"
$t.Cell(18,1).Range.Text = "75-73=2"
$t.Cell(18,2).Range.Text = "36+48=84"
$t.Cell(18,3).Range.Text = "21+8=29"
$t.Cell(18,4).Range.Text = "78-59=19"
$t.Cell(18,5).Range.Text = "73-5=68"
$t.Cell(19,1).Range.Text = "67-30=37"
$t.Cell(19,2).Range.Text = "21+11=32"
$t.Cell(19,3).Range.Text = "75-34=41"
$t.Cell(19,4).Range.Text = "83-43=40"
$t.Cell(19,5).Range.Text = "98-39=59"
$t.Cell(20,1).Range.Text = "45+5=50"
$t.Cell(20,2).Range.Text = "82-42=40"
$t.Cell(20,3).Range.Text = "3+50=53"
$t.Cell(20,4).Range.Text = "83-32=51"
$t.Cell(20,5).Range.Text = "34+41=75"
